# DB Backup - Stake and cus position
# Update the stakeholder member code and customer id values, then leave the
# selection/active sheet the way the author left it before saving.

$wb = $excel.ActiveWorkbook

# project_stakeholder: stakeholder_member_code value changed
$wsStakeholder = $wb.Worksheets.Item("project_stakeholder")
$wsStakeholder.Range("A2").Value = "2020-C7DVWDZ7"

# customer: customer_id values changed (row 2 and row 3)
$wsCustomer = $wb.Worksheets.Item("customer")
$wsCustomer.Range("A2").Value = "2019-UIDMS"
$wsCustomer.Range("A3").Value = "2020-JH9UHG26"

# Leave the cursor parked on project_stakeholder!C6 ...
$wsStakeholder.Activate() | Out-Null
$wsStakeholder.Range("C6").Select() | Out-Null

# ... then move on to customer!E6, which ends up the active sheet/selection
# at save time.
$wsCustomer.Activate() | Out-Null
$wsCustomer.Range("E6").Select() | Out-Null
